$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(147, 1).Value = 147
$ws.Cells.Item(147, 2).Value = "Biobío"
$ws.Cells.Item(147, 3).Value = "Temas de Personas"
$ws.Cells.Item(147, 4).Value = "Procesos de selección: se están cerrando los procesos de entrevistas de DAF, el de plataforma está en evaluación cv."
$ws.Cells.Item(147, 5).Value = "Pendiente"
$ws.Cells.Item(147, 6).Value = 0
$ws.Cells.Item(147, 7).Value = "15-09-2025"

$ws.Cells.Item(148, 1).Value = 148
$ws.Cells.Item(148, 2).Value = "Biobío"
$ws.Cells.Item(148, 3).Value = "Ejecución Presupuestaria"
$ws.Cells.Item(148, 4).Value = "Licitación privada para compras de capacitaciones: en acta de evaluación. Línea 1 y 2 con listado de inscritos. Y el resto 1 capacitación semanal para los servicios de salud. `nLicitación de cursos idea de realizarla por dos años. Más parrilla de cursos a demanda, se ampliará la oferta. Pausas saludables: en licitación son 65 millones. Se va a analizar cómo se está usando. SLEP de Puelche, se está evaluando si se va a quedar o no. Con el SLEP los Copihues, se está ofreciendo la máxima colaboración."
$ws.Cells.Item(148, 5).Value = "Pendiente"
$ws.Cells.Item(148, 6).Value = 0
$ws.Cells.Item(148, 7).Value = "15-09-2025"

$ws.Cells.Item(149, 1).Value = 149
$ws.Cells.Item(149, 2).Value = "Biobío"
$ws.Cells.Item(149, 3).Value = "Otros"
$ws.Cells.Item(149, 4).Value = "Temas exploración prestadores públicos: traumatológico y COSAM (hay 3 COSAM en Concepción). La canasta de prestaciones traumatológicas opera con arancel particular. Quedó en el convenio original con el Ministerio de Salud. Insumos y medicamentos a valor particular también."
$ws.Cells.Item(149, 5).Value = "Pendiente"
$ws.Cells.Item(149, 6).Value = 0
$ws.Cells.Item(149, 7).Value = "15-09-2025"

$ws.Cells.Item(150, 1).Value = 150
$ws.Cells.Item(150, 2).Value = "Biobío"
$ws.Cells.Item(150, 3).Value = "Infraestructura"
$ws.Cells.Item(150, 4).Value = "Podría la persona que está pidiendo traslado de Los Ríos, ejecutar las verificaciones de medidas (salud). Los Ángeles es donde se requiere apoyo de prevención (se debería solicitar cupo). Con la nueva licitación de psicólogos se requiere en Los Ángeles)."
$ws.Cells.Item(150, 5).Value = "Pendiente"
$ws.Cells.Item(150, 6).Value = 0
$ws.Cells.Item(150, 7).Value = "15-09-2025"

$ws.Cells.Item(151, 1).Value = 151
$ws.Cells.Item(151, 2).Value = "Biobío"
$ws.Cells.Item(151, 3).Value = "Infraestructura"
$ws.Cells.Item(151, 4).Value = "Se ha ordenado bodegas y se eliminaron los documentos no necesarios. Se ha realizado ciertas mejoras de contar con techos en salidas.`nContar con sucursal en Talcahuano, para potenciar ámbito preventivo en la zona, con punto de atención."
$ws.Cells.Item(151, 5).Value = "Pendiente"
$ws.Cells.Item(151, 6).Value = 0
$ws.Cells.Item(151, 7).Value = "15-09-2025"

$ws.Cells.Item(152, 1).Value = 152
$ws.Cells.Item(152, 2).Value = "Los Ríos"
$ws.Cells.Item(152, 3).Value = "Temas de Personas"
$ws.Cells.Item(152, 4).Value = "Se debe resolver tema movilidades internas de Cristian Herrera y Lissette Latorre, directora enviará correo a Natalia."
$ws.Cells.Item(152, 5).Value = "Pendiente"
$ws.Cells.Item(152, 6).Value = 0
$ws.Cells.Item(152, 7).Value = "15-09-2025"

$ws.Cells.Item(153, 1).Value = 153
$ws.Cells.Item(153, 2).Value = "Los Ríos"
$ws.Cells.Item(153, 3).Value = "Indicadores de desempeño"
$ws.Cells.Item(153, 4).Value = "Respecto del plan de prevención: hay un tema de retroalimentación desde NC de las planillas de los indicadores psicosocial. También con tema servicios públicos, hay una actualización que realizó la región y el indicador aumentó. Sílice está bajo para que puedan rastrear nuevamente."
$ws.Cells.Item(153, 5).Value = "Pendiente"
$ws.Cells.Item(153, 6).Value = 0
$ws.Cells.Item(153, 7).Value = "15-09-2025"

$ws.Cells.Item(154, 1).Value = 154
$ws.Cells.Item(154, 2).Value = "Los Ríos"
$ws.Cells.Item(154, 3).Value = "Temas de Personas"
$ws.Cells.Item(154, 4).Value = "Resolver jefatura sección de prevención."
$ws.Cells.Item(154, 5).Value = "Pendiente"
$ws.Cells.Item(154, 6).Value = 0
$ws.Cells.Item(154, 7).Value = "15-09-2025"

$ws.Cells.Item(155, 1).Value = 155
$ws.Cells.Item(155, 2).Value = "Los Ríos"
$ws.Cells.Item(155, 3).Value = "Político Institucional"
$ws.Cells.Item(155, 4).Value = "Caso Pablo Ulloa: Lo fueron a visitar para actualizar el informe social. Lo pendiente era su internet, el que se debe reembolsar de acuerdo a instrucciones, pero es una antena satelital que se instaló el año pasado. Se pide a la directora regional escalar con salud, enviará correo."
$ws.Cells.Item(155, 5).Value = "Pendiente"
$ws.Cells.Item(155, 6).Value = 0
$ws.Cells.Item(155, 7).Value = "15-09-2025"
